$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Matrix area (K7:S15): re-enter the counting formulas so the K column
# (K8:K15) becomes one contiguous shared formula down the column, matching
# a "fill down" of K8's formula through K15. ---
$ws.Range("K8:K15").Formula = "=S7+1"

$ws.Range("L8").Formula = "=K8+1"
$ws.Range("M8:S8").Formula = "=L8+1"

$ws.Range("L9").Formula = "=K9+1"
$ws.Range("M9:S9").Formula = "=L9+1"

$ws.Range("L10:Q10").Formula = "=K10+1"
$ws.Range("R10").Formula = "=Q10+1"
$ws.Range("S10").Formula = "=R10+1"
$ws.Range("L11:S13").Formula = "=K11+1"

$ws.Range("L14:S14").Formula = "=K14+1"
$ws.Range("L15:S15").Formula = "=K15+1"

# --- "Convert coordinate to node" calculator block (rows 18-21) ---
# Existing x/y inputs reset to 0
$ws.Range("C18").Value = 0
$ws.Range("C19").Value = 0

# New "node -> coordinate" helper block in K:O
# Row 19 first so the new shared string "y2 =" is interned before "x2 ="
$ws.Range("K19").Value = "y="
$ws.Range("L19").Value = 14
$ws.Range("N19").Value = "y2 ="
$ws.Range("O19").Formula = "=L19/2"

$ws.Range("K18").Value = "x="
$ws.Range("L18").Value = 16
$ws.Range("N18").Value = "x2 ="
$ws.Range("O18").Formula = "=L18/2"

# Replace the old literal "2,3" label cell with the node-number formula
$ws.Range("J21").Formula = "=O19*9+O18"

# --- Selection / view state ---
$ws.Range("T8").Select()
